# Update cryptos.xlsx cell values per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.643.02'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '2.098.47'
$ws.Range('E3').Value = '  +10.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.668'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.81'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '60.38'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.369'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0729'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.89%  '
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '2.405.43'
$ws.Range('E14').Value = '  +11.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.835'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.81%  '
$ws.Range('D16').Value = '2.095.94'
$ws.Range('E16').Value = '  +10.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '36.742.51'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('E25').Value = '  -10.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.65%  '
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.49%  '
$ws.Range('E30').Value = '  -4.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.04'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +55.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0593'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0894'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.89'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.36%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +20.09%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.908'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.98%  '
$ws.Range('E39').Value = '  -5.90%  '
$ws.Range('E40').Value = '  -9.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0218'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +16.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.11'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.60%  '
$ws.Range('D46').Value = '1.367.25'
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.291.67'
$ws.Range('E48').Value = '  +10.80%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +16.00%  '
